# Update the LCFC scouting report with the new player's details
# (fix summary average value, per commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / fixture info ---
$ws.Range("H1").Value = "Chelsea vs Arsenal"      # Fixture
$ws.Range("G4").Value = "Sean Raisi"              # Name
$ws.Range("L4").Value = "'19"                     # Age (kept as text)
$ws.Range("G5").Value = "'5.7"                    # Height (kept as text)
$ws.Range("L5").Value = "Chelsea"                 # Club
$ws.Range("L6").Value = "2-0"                     # H/T
$ws.Range("G7").Value = "Arsenal"                 # Playing Against
$ws.Range("L7").Value = "3-0"                     # F/T
$ws.Range("G8").Value = "22/01/18"                # Date

# --- Attribute score grid: every scored cell becomes 9 (kept as text) ---
$scoreCells = @(
    "C13","F13","I13","L13","O13","R13",
    "C14","F14","I14","L14","O14","R14",
    "C15","F15","I15","L15","O15","R15",
    "C16","F16","I16","L16","O16","R16",
    "C17","F17","L17","O17","R17",
    "F18","R18"
)
foreach ($ref in $scoreCells) {
    $ws.Range($ref).Value = "'9"
}

# --- Notes section ---
$ws.Range("A23").Value = ""

# --- Summary text ---
$ws.Range("A27").Value = "Raisi, Sean was scouted playing for Chelsea on 22/01/18. Raisi, Sean performed to grade A with an average score of 9 showing some outstanding attributes."

# --- Player rating ---
$ws.Range("H30").Value = "A"
